$d = $word.ActiveDocument

# 1. Title heading + bold run near the end (replace all occurrences)
$d.Content.Find.Execute("Play Jaguar Gold for Free - Review of Skywind Slot", $true, $false, $false, $false, $false, $true, 1, $false, "Play Jaguar Gold Free: Review of Gameplay, Graphics, and Wins", 2) | Out-Null

# 2. "What we like" bullet 1
$d.Content.Find.Execute("Stunning graphics and engaging jungle theme", $true, $false, $false, $false, $false, $true, 1, $false, "Stunning graphics and engaging gameplay", 2) | Out-Null

# 3. "What we like" bullet 3
$d.Content.Find.Execute("Jumbo Links Bonus and Diamond Bonus offer exciting opportunities", $true, $false, $false, $false, $false, $true, 1, $false, "Free spins and Diamond Bonus add excitement", 2) | Out-Null

# 4. "What we like" bullet 4
$d.Content.Find.Execute("Suitable for both experienced and less experienced players", $true, $false, $false, $false, $false, $true, 1, $false, "Jumbo Links Bonus and multipliers for big wins", 2) | Out-Null

# 5. "What we don't like" bullet 2
$d.Content.Find.Execute("Limited number of Bonus symbols to trigger Jumbo Links Bonus", $true, $false, $false, $false, $false, $true, 1, $false, "Requires patience to achieve big wins", 2) | Out-Null

# 6. Italic summary paragraph
$d.Content.Find.Execute("Discover the jungle adventure with Jaguar Gold. Play for free and enjoy special functions, bonuses, and stunning graphics in this Skywind Group slot game.", $true, $false, $false, $false, $false, $true, 1, $false, "Discover the exciting gameplay, stunning graphics, and winning potential of Jaguar Gold. Play free now!", 2) | Out-Null
